# Generate Report for Handoff
#
# Updates the "Latest Handoff Date" / "Latest Handoff Datetime" timestamps for the
# 5c4cc5a0-b7f6-4851-8ce5-f381df8f46f4.md file, which has just been handed off again,
# on the Overview sheet as well as the per-locale (zh-cn, de-de) detail sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 5 is the 5c4cc5a0-b7f6-4851-8ce5-f381df8f46f4.md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D5").Value = "2016-27-12 02:27:13"

# --- zh-cn sheet: row 5 is the 5c4cc5a0-... source file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E5").Value = "2016-03-12 02:27:10"

# --- de-de sheet: row 5 is the 5c4cc5a0-... source file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E5").Value = "2016-03-12 02:27:13"
